# Fruta / hortaliza, semanal
# The weekly refresh reshuffles the per-row data (Fecha, Volumen, Precio
# minimo/maximo/promedio, Origen, Precio $/Kg) among the existing data rows
# (rows 2-19) of the sheet. Columns A,B,C,E,F,G,H,I,J,K,L,Q,T are identical
# for every row and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D, M, N, O, P, R, S
$rows = @{
    2  = @{ D = 44614; M = 45;  N = 6000; O = 6000; P = 6000; R = "Provincia de Linares"; S = 3000 }
    3  = @{ D = 44627; M = 45;  N = 6000; O = 6000; P = 6000; R = "Provincia de Linares"; S = 3000 }
    4  = @{ D = 45001; M = 66;  N = 7500; O = 8000; P = 7773; R = "Provincia de Curicó";  S = 3886 }
    5  = @{ D = 44214; M = 48;  N = 6000; O = 6000; P = 6000; R = "Provincia de Linares"; S = 3000 }
    6  = @{ D = 44586; M = 80;  N = 7000; O = 7000; P = 7000; R = "Provincia de Curicó";  S = 3500 }
    7  = @{ D = 44959; M = 40;  N = 7000; O = 7000; P = 7000; R = "Provincia de Curicó";  S = 3500 }
    8  = @{ D = 44588; M = 160; N = 6500; O = 7000; P = 6750; R = "Provincia de Curicó";  S = 3375 }
    9  = @{ D = 44628; M = 40;  N = 6000; O = 6000; P = 6000; R = "Provincia de Linares"; S = 3000 }
    10 = @{ D = 44587; M = 165; N = 6500; O = 7000; P = 6742; R = "Provincia de Linares"; S = 3371 }
    11 = @{ D = 44589; M = 60;  N = 6000; O = 6000; P = 6000; R = "Provincia de Curicó";  S = 3000 }
    12 = @{ D = 44606; M = 45;  N = 7000; O = 7000; P = 7000; R = "Provincia de Linares"; S = 3500 }
    13 = @{ D = 44582; M = 150; N = 6000; O = 6500; P = 6233; R = "Provincia de Curicó";  S = 3116 }
    14 = @{ D = 44209; M = 58;  N = 6000; O = 6000; P = 6000; R = "Provincia de Curicó";  S = 3000 }
    15 = @{ D = 44960; M = 40;  N = 7000; O = 7000; P = 7000; R = "Provincia de Curicó";  S = 3500 }
    16 = @{ D = 44592; M = 30;  N = 8000; O = 8000; P = 8000; R = "Provincia de Linares"; S = 4000 }
    17 = @{ D = 44211; M = 45;  N = 6000; O = 6000; P = 6000; R = "Provincia de Curicó";  S = 3000 }
    18 = @{ D = 44585; M = 160; N = 6500; O = 7000; P = 6750; R = "Provincia de Curicó";  S = 3375 }
    19 = @{ D = 44974; M = 130; N = 7000; O = 7500; P = 7269; R = "Provincia de Curicó";  S = 3634 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.S
}
